$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '26.574.46'
$ws.Cells.Item(2, 4).NumberFormat = "General"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '  -8.44%  '
$ws.Cells.Item(2, 5).NumberFormat = "General"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.654.38'
$ws.Cells.Item(3, 4).NumberFormat = "General"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '  -8.99%  '
$ws.Cells.Item(3, 5).NumberFormat = "General"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.008'
$ws.Cells.Item(4, 4).NumberFormat = "General"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '  +0.34%  '
$ws.Cells.Item(4, 5).NumberFormat = "General"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '219.91'
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '  -5.50%  '
$ws.Cells.Item(5, 5).NumberFormat = "General"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.5050'
$ws.Cells.Item(6, 4).NumberFormat = "General"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '  -14.45%  '
$ws.Cells.Item(6, 5).NumberFormat = "General"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.008'
$ws.Cells.Item(7, 4).NumberFormat = "General"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '  +0.38%  '
$ws.Cells.Item(7, 5).NumberFormat = "General"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2551'
$ws.Cells.Item(8, 4).NumberFormat = "General"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '  -7.10%  '
$ws.Cells.Item(8, 5).NumberFormat = "General"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '21.58'
$ws.Cells.Item(9, 4).NumberFormat = "General"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '  -5.87%  '
$ws.Cells.Item(9, 5).NumberFormat = "General"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.06131'
$ws.Cells.Item(10, 4).NumberFormat = "General"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '  -9.36%  '
$ws.Cells.Item(10, 5).NumberFormat = "General"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07339'
$ws.Cells.Item(11, 4).NumberFormat = "General"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '  -2.07%  '
$ws.Cells.Item(11, 5).NumberFormat = "General"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.651.16'
$ws.Cells.Item(12, 4).NumberFormat = "General"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '  -8.98%  '
$ws.Cells.Item(12, 5).NumberFormat = "General"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.433'
$ws.Cells.Item(13, 4).NumberFormat = "General"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '  -5.09%  '
$ws.Cells.Item(13, 5).NumberFormat = "General"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.5730'
$ws.Cells.Item(14, 4).NumberFormat = "General"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '  -7.84%  '
$ws.Cells.Item(14, 5).NumberFormat = "General"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '1.876.36'
$ws.Cells.Item(15, 4).NumberFormat = "General"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '  -9.12%  '
$ws.Cells.Item(15, 5).NumberFormat = "General"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.000007981'
$ws.Cells.Item(16, 4).NumberFormat = "General"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '  -15.04%  '
$ws.Cells.Item(16, 5).NumberFormat = "General"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '64.30'
$ws.Cells.Item(17, 4).NumberFormat = "General"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '  -13.68%  '
$ws.Cells.Item(17, 5).NumberFormat = "General"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '26.581.16'
$ws.Cells.Item(18, 4).NumberFormat = "General"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '  -7.59%  '
$ws.Cells.Item(18, 5).NumberFormat = "General"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.917'
$ws.Cells.Item(19, 4).NumberFormat = "General"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '  -9.32%  '
$ws.Cells.Item(19, 5).NumberFormat = "General"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '1.009'
$ws.Cells.Item(20, 4).NumberFormat = "General"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '  +0.47%  '
$ws.Cells.Item(20, 5).NumberFormat = "General"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '  -7.26%  '
$ws.Cells.Item(21, 5).NumberFormat = "General"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '178.95'
$ws.Cells.Item(22, 4).NumberFormat = "General"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '  -13.75%  '
$ws.Cells.Item(22, 5).NumberFormat = "General"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.009'
$ws.Cells.Item(23, 4).NumberFormat = "General"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '  +0.34%  '
$ws.Cells.Item(23, 5).NumberFormat = "General"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.192'
$ws.Cells.Item(24, 4).NumberFormat = "General"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '  -8.32%  '
$ws.Cells.Item(24, 5).NumberFormat = "General"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '142.65'
$ws.Cells.Item(25, 4).NumberFormat = "General"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '  -7.99%  '
$ws.Cells.Item(25, 5).NumberFormat = "General"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '7.515'
$ws.Cells.Item(26, 4).NumberFormat = "General"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '  -3.20%  '
$ws.Cells.Item(26, 5).NumberFormat = "General"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.1145'
$ws.Cells.Item(27, 4).NumberFormat = "General"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '  -9.71%  '
$ws.Cells.Item(27, 5).NumberFormat = "General"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '14.91'
$ws.Cells.Item(28, 4).NumberFormat = "General"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '  -8.24%  '
$ws.Cells.Item(28, 5).NumberFormat = "General"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.324'
$ws.Cells.Item(29, 4).NumberFormat = "General"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '  -5.50%  '
$ws.Cells.Item(29, 5).NumberFormat = "General"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.05803'
$ws.Cells.Item(30, 4).NumberFormat = "General"
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = '  -9.84%  '
$ws.Cells.Item(30, 5).NumberFormat = "General"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.337'
$ws.Cells.Item(31, 4).NumberFormat = "General"
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = '  -6.46%  '
$ws.Cells.Item(31, 5).NumberFormat = "General"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.405'
$ws.Cells.Item(32, 4).NumberFormat = "General"
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = '  -8.27%  '
$ws.Cells.Item(32, 5).NumberFormat = "General"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.371'
$ws.Cells.Item(33, 4).NumberFormat = "General"
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = '  -8.09%  '
$ws.Cells.Item(33, 5).NumberFormat = "General"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.581'
$ws.Cells.Item(34, 4).NumberFormat = "General"
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = '  -5.71%  '
$ws.Cells.Item(34, 5).NumberFormat = "General"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.9728'
$ws.Cells.Item(35, 4).NumberFormat = "General"
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = '  -7.06%  '
$ws.Cells.Item(35, 5).NumberFormat = "General"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.425'
$ws.Cells.Item(36, 4).NumberFormat = "General"
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = '  -4.39%  '
$ws.Cells.Item(36, 5).NumberFormat = "General"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.5947'
$ws.Cells.Item(37, 4).NumberFormat = "General"
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '  -5.33%  '
$ws.Cells.Item(37, 5).NumberFormat = "General"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.634'
$ws.Cells.Item(38, 4).NumberFormat = "General"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01585'
$ws.Cells.Item(39, 4).NumberFormat = "General"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '  -6.52%  '
$ws.Cells.Item(39, 5).NumberFormat = "General"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.068.05'
$ws.Cells.Item(40, 4).NumberFormat = "General"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '  -5.27%  '
$ws.Cells.Item(40, 5).NumberFormat = "General"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.8598'
$ws.Cells.Item(41, 4).NumberFormat = "General"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '  -0.95%  '
$ws.Cells.Item(41, 5).NumberFormat = "General"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.010'
$ws.Cells.Item(42, 4).NumberFormat = "General"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '  +0.63%  '
$ws.Cells.Item(42, 5).NumberFormat = "General"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '5.695'
$ws.Cells.Item(43, 4).NumberFormat = "General"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '  -11.13%  '
$ws.Cells.Item(43, 5).NumberFormat = "General"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '96.64'
$ws.Cells.Item(44, 4).NumberFormat = "General"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '  -3.39%  '
$ws.Cells.Item(44, 5).NumberFormat = "General"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.791.36'
$ws.Cells.Item(45, 4).NumberFormat = "General"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '  -9.25%  '
$ws.Cells.Item(45, 5).NumberFormat = "General"
$ws.Cells.Item(46, 2).NumberFormat = "@"
$ws.Cells.Item(46, 2).Value = 'Frax'
$ws.Cells.Item(46, 2).NumberFormat = "General"
$ws.Cells.Item(46, 3).NumberFormat = "@"
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(46, 3).NumberFormat = "General"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.010'
$ws.Cells.Item(46, 4).NumberFormat = "General"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '  +0.80%  '
$ws.Cells.Item(46, 5).NumberFormat = "General"
$ws.Cells.Item(47, 2).NumberFormat = "@"
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 2).NumberFormat = "General"
$ws.Cells.Item(47, 3).NumberFormat = "@"
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 3).NumberFormat = "General"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '55.47'
$ws.Cells.Item(47, 4).NumberFormat = "General"
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '  -7.64%  '
$ws.Cells.Item(47, 5).NumberFormat = "General"
$ws.Cells.Item(48, 2).NumberFormat = "@"
$ws.Cells.Item(48, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(48, 2).NumberFormat = "General"
$ws.Cells.Item(48, 3).NumberFormat = "@"
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(48, 3).NumberFormat = "General"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.00000000105'
$ws.Cells.Item(48, 4).NumberFormat = "General"
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '  -7.13%  '
$ws.Cells.Item(48, 5).NumberFormat = "General"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.4371'
$ws.Cells.Item(49, 4).NumberFormat = "General"
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '  -3.10%  '
$ws.Cells.Item(49, 5).NumberFormat = "General"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.05188'
$ws.Cells.Item(50, 4).NumberFormat = "General"
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '  -5.08%  '
$ws.Cells.Item(50, 5).NumberFormat = "General"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '7.727'
$ws.Cells.Item(51, 4).NumberFormat = "General"
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '  -6.18%  '
$ws.Cells.Item(51, 5).NumberFormat = "General"
